# "Generate Report for handback"
#
# This script records a handback for the single target file (a.md.md) in
# both locale sheets (zh-cn, de-de):
#   - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   - "Latest Target File" (E) / "Latest Handback File" (F) are populated with
#     the handed-back file name / xlf name (mirroring the handoff columns,
#     with their own hyperlinks)
#   - "Latest Handback DateTime" (G) is stamped with the actual handback time

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# The "Overview" rollup sheet mirrors each locale's Status column (it shares
# the same "Ready for handoff" text) - update it in lockstep so no sheet is
# left displaying the stale status.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

function Update-LocaleSheet($SheetName, $XlfName, $HandbackTime, $AMdTargetUrl, $XlfTargetUrl) {
    $ws = $wb.Worksheets.Item($SheetName)

    # Status column (B) for the two real rows: row2 = a.md.md, row3 = b.md.md
    $ws.Range("B2").Value = $newStatus
    $ws.Range("B3").Value = $newStatus

    # Latest Target File (E) / Latest Handback File (F) for row 2
    $ws.Range("E2").Value = "a.md.md"
    $ws.Range("F2").Value = $XlfName
    $ws.Hyperlinks.Add($ws.Range("E2"), $AMdTargetUrl, "", "", "a.md.md")
    $ws.Hyperlinks.Add($ws.Range("F2"), $XlfTargetUrl, "", "", $XlfName)

    # Latest Target File (E) / Latest Handback File (F) for row 3
    $ws.Range("E3").Value = "a.md.md"
    $ws.Range("F3").Value = $XlfName
    $ws.Hyperlinks.Add($ws.Range("E3"), $AMdTargetUrl, "", "", "a.md.md")
    $ws.Hyperlinks.Add($ws.Range("F3"), $XlfTargetUrl, "", "", $XlfName)

    # Latest Handback DateTime (G) for rows 2 and 3
    $ws.Range("G2").Value = $HandbackTime
    $ws.Range("G3").Value = $HandbackTime
}

Update-LocaleSheet "zh-cn" "a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf" "2016-01-14 03:09:38" "https://github.com/OpenLocalizationTest/oltest/blob/cafddb680a83aa4d7bfb0993a974de43ae9670ea/e2e/a.md.md" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f7103e36fcf5b5328e0c111f6873fabb13cb981a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/a.md.370104d57010292b5765347db07350cde3a977e8.zh-cn.xlf"

Update-LocaleSheet "de-de" "a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf" "2016-01-14 03:09:58" "https://github.com/OpenLocalizationTest/oltest/blob/cafddb680a83aa4d7bfb0993a974de43ae9670ea/e2e/a.md.md" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d7a3bbd0db394824eb9be2a98b5e5f32eea1ad36/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/a.md.370104d57010292b5765347db07350cde3a977e8.de-de.xlf"
